$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.067.32'
$ws.Range('E2').Value = '  +0.33%  '

# Row 3
$ws.Range('D3').Value = '3.399.10'
$ws.Range('E3').Value = '  -2.57%  '

# Row 4
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.96'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.23%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.63'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.14%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.623'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +4.88%  '

# Row 8
$ws.Range('E8').Value = '  -0.04%  '

# Row 9
$ws.Range('D9').Value = '3.396.74'
$ws.Range('E9').Value = '  -2.53%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.131'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.48%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.94'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.29%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.414'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.67%  '

# Row 13
$ws.Range('D13').Value = '3.986.18'
$ws.Range('E13').Value = '  -2.67%  '

# Row 14
$ws.Range('E14').Value = '  +0.75%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.05'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.79%  '

# Row 16
$ws.Range('D16').Value = '66.124.56'
$ws.Range('E16').Value = '  +0.22%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000172'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.36%  '

# Row 18
$ws.Range('D18').Value = '3.394.07'
$ws.Range('E18').Value = '  -2.75%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.88'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.75%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.66'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.51%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '366.52'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.16%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.55'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.76%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.88'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.10%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.10%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.529'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.82%  '

# Row 26
$ws.Range('E26').Value = '  +1.68%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.75'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.60%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.180'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.20%  '

# Row 29
$ws.Range('E29').Value = '  +0.03%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.99'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.22%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.73'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.72%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '23.17'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.15%  '

# Row 33
$ws.Range('E33').Value = '  +0.06%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.99'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.91%  '

# Row 35
$ws.Range('E35').Value = '  -3.01%  '

# Row 36
$ws.Range('E36').Value = '  -2.03%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '161.46'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.47%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.857'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.35%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.17'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -8.02%  '

# Row 40
$ws.Range('E40').Value = '  -0.37%  '

# Row 41
$ws.Range('E41').Value = '  +2.08%  '

# Row 42
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.680.29'
$ws.Range('E42').Value = '  -3.89%  '

# Row 43
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.35'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.30%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.20'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.85%  '

# Row 45
$ws.Range('E45').Value = '  -1.41%  '

# Row 46
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '336.53'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +9.90%  '

# Row 47
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '39.72'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.34%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '24.53'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.26%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0285'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.16%  '

# Row 50
$ws.Range('E50').Value = '  +2.87%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '31.39'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.89%  '
